$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell rewrites (row => new text)
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "106"

$t.Cell(6,1).Range.Text  = "0.00015"
$t.Cell(7,1).Range.Text  = "0.00006"
$t.Cell(9,1).Range.Text  = "0.00004"
$t.Cell(10,1).Range.Text = ""
$t.Cell(11,1).Range.Text = "0.00015"
$t.Cell(12,1).Range.Text = "0.00484"

# Collapse the three wide tab-separated summary rows down to a single value
$t.Cell(44,1).Range.Text = "99.99"
$t.Cell(45,1).Range.Text = "0"
$t.Cell(46,1).Range.Text = "92"
